$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts the assigned string
# into a numeric value (losing the literal formatting, e.g. "1.00" -> 1).
$ws.Range("D2").Value = "59.247.47"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.602.98"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.84"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.57"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "3.061.07"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "59.187.37"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "2.614.73"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.38"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.13"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.39"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.61"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.20"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("E28").Value = "  +4.12%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +7.66%  "
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.73"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.82"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.99"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.834"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.827"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.57"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "275.05"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("D47").Value = "1.954.83"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0224"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.55"
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.52"
$ws.Range("E51").Value = "  -0.13%  "
